# Adapt column header formatting to respective input file names.
# - Rename the "_old" / "_new" header-name suffixes to the concrete
#   format-version identifiers "_FV2410" / "_FV2504".
# - Freeze the header row (row 1) in the sheet view.
# - Wrap the data range in an Excel Table ("Table1") so the (renamed)
#   header labels become the table's column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row A1:U1 -------------------------------------------
# Ordered (address, new header text) pairs -- avoids relying on hashtable
# enumeration order, which PowerShell does not guarantee.
$headerRenames = @(
    ,@("A1", "Segmentname_FV2410")
    ,@("B1", "Segmentgruppe_FV2410")
    ,@("C1", "Segment_FV2410")
    ,@("D1", "Datenelement_FV2410")
    ,@("E1", "Segment ID_FV2410")
    ,@("F1", "Code_FV2410")
    ,@("G1", "Qualifier_FV2410")
    ,@("H1", "Beschreibung_FV2410")
    ,@("I1", "Bedingungsausdruck_FV2410")
    ,@("J1", "Bedingung_FV2410")
    ,@("K1", "diff")
    ,@("L1", "Segmentname_FV2504")
    ,@("M1", "Segmentgruppe_FV2504")
    ,@("N1", "Segment_FV2504")
    ,@("O1", "Datenelement_FV2504")
    ,@("P1", "Segment ID_FV2504")
    ,@("Q1", "Code_FV2504")
    ,@("R1", "Qualifier_FV2504")
    ,@("S1", "Beschreibung_FV2504")
    ,@("T1", "Bedingungsausdruck_FV2504")
    ,@("U1", "Bedingung_FV2504")
)

foreach ($pair in $headerRenames) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- 2) Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Turn the used range into an Excel Table -----------------------------
$tableRange = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

$ws.Range("A1").Select()
